# Update countries & provincias Spain
# Refresh the COVID dashboard data on sheet "Pais": several countries swap
# rank (their row order), several rows get refreshed case counts, and the
# "last updated" timestamp is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Kazajistan (row 34) - refreshed totals
$ws.Range("B34").Value = 107199
$ws.Range("C34").Value = 65
$ws.Range("D34").Value = 101822
$ws.Range("E34").Value = 3706

# Belgica overtakes Kuwait -> rows 38/39 swap label + stats
$ws.Range("A38").Value = "Belgica"
$ws.Range("B38").Value = 99649
$ws.Range("C38").Value = 1673
$ws.Range("D38").Value = 18908
$ws.Range("E38").Value = 70804
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 9937

$ws.Range("A39").Value = "Kuwait"
$ws.Range("B39").Value = 98528
$ws.Range("D39").Value = 88776
$ws.Range("E39").Value = 9172
$ws.Range("H39").Value = 580

# Honduras (row 50) - refreshed totals
$ws.Range("B50").Value = 70611
$ws.Range("C50").Value = 491
$ws.Range("D50").Value = 21149
$ws.Range("E50").Value = 47316
$ws.Range("G50").Value = 24
$ws.Range("H50").Value = 2146

# Venezuela (row 53) - refreshed totals
$ws.Range("B53").Value = 65174
$ws.Range("D53").Value = 54218
$ws.Range("E53").Value = 10426
$ws.Range("H53").Value = 530

# Birmania overtakes Jamaica -> rows 125/126 swap label + stats
$ws.Range("A125").Value = "Birmania"
$ws.Range("B125").Value = 4621
$ws.Range("C125").Value = 154
$ws.Range("D125").Value = 1130
$ws.Range("E125").Value = 3416
$ws.Range("G125").Value = 5
$ws.Range("H125").Value = 75

$ws.Range("A126").Value = "Jamaica"
$ws.Range("B126").Value = 4571
$ws.Range("D126").Value = 1264
$ws.Range("E126").Value = 3252
$ws.Range("H126").Value = 55

# Belice (row 158) - refreshed totals
$ws.Range("B158").Value = 1590
$ws.Range("C158").Value = 23
$ws.Range("D158").Value = 812
$ws.Range("E158").Value = 758
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 20

# Islas Turcas y Caicos (row 172) - refreshed totals
$ws.Range("B172").Value = 667
$ws.Range("C172").Value = 4
$ws.Range("D172").Value = 567
$ws.Range("E172").Value = 95

# San Martin (Parte Holandesa) (row 173) - refreshed totals
$ws.Range("B173").Value = 574
$ws.Range("C173").Value = 9
$ws.Range("D173").Value = 488
$ws.Range("E173").Value = 66
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = 20

# Santa Lucia / Timor Oriental swap places (rows 204/205), stats unchanged
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# Montserrat overtakes Islas Malvinas -> rows 214/215 swap label + stats
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Bump the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 05:28"
